# PSY_age_of_onset.xlsx — rename "Schizotypic" row label to "Schizotypy"
# and update the saved sheet view (scroll position + active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix -----------------------------------------------------
# Row 11 currently reads "Schizotypic"; the correct term is "Schizotypy".
# (Row 12, "PTSD", is untouched — shared-string bookkeeping is handled
# automatically by the engine when the old string becomes unused.)
$ws.Range("A11").Value = "Schizotypy"

# --- View state --------------------------------------------------------
# Scroll the sheet so row 5 is at the top and select H10, matching the
# author's last on-screen state when the workbook was saved.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H10").Select()
